# Resultados_SmartScore.xlsx -- "Actualizacion SmartScore desde Streamlit (Miranda)"
#
# Miranda resubmitted the SmartScore form. Her previous response is preserved
# as a new row 7 (same answers, original text-typed SmartScore formatting, new
# submission timestamp), and row 6 is corrected in place so its SmartScore
# columns (G, J, M, P, S, V, Y, AB, AE) are stored as real numbers instead of
# text, matching every other row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: append new submission row 7 -- identical answers to row 6, with
#             the new timestamp. SmartScore columns are written as literal
#             text ("0.510", not 0.51) to match how the form originally saved
#             them, before the row-6 type fix below. ---
$ws.Cells.Item(7, 1).Value = "Miranda"
$ws.Cells.Item(7, 2).Value = 25
$ws.Cells.Item(7, 3).Value = "Femenino"
$ws.Cells.Item(7, 4).Value = "2025-10-28 05:58:21"
$ws.Cells.Item(7, 5).Value = "{
  `"portion`": 0.8,
  `"diet`": 0.5714285714285714,
  `"salt`": 0.6,
  `"fat`": 0.8,
  `"natural`": 0.6,
  `"convenience`": 0.4,
  `"price`": 0.8
}"
$ws.Cells.Item(7, 6).Value = "Nongshim Neoguri Spicy Seafood"
$ws.Cells.Item(7, 7).Value = "'0.575"
$ws.Cells.Item(7, 7).Style = "Normal"
$ws.Cells.Item(7, 8).Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
$ws.Cells.Item(7, 9).Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Cells.Item(7, 10).Value = "'0.510"
$ws.Cells.Item(7, 10).Style = "Normal"
$ws.Cells.Item(7, 11).Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"
$ws.Cells.Item(7, 12).Value = "Maruchan Ramen Sabor Pollo"
$ws.Cells.Item(7, 13).Value = "'0.509"
$ws.Cells.Item(7, 13).Style = "Normal"
$ws.Cells.Item(7, 14).Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"
$ws.Cells.Item(7, 15).Value = "Kraft Macaroni & Cheese Dinner"
$ws.Cells.Item(7, 16).Value = "'0.650"
$ws.Cells.Item(7, 16).Style = "Normal"
$ws.Cells.Item(7, 17).Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
$ws.Cells.Item(7, 18).Value = "Annie’s Shells & White Cheddar"
$ws.Cells.Item(7, 19).Value = "'0.587"
$ws.Cells.Item(7, 19).Style = "Normal"
$ws.Cells.Item(7, 20).Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"
$ws.Cells.Item(7, 21).Value = "Amy’s Macaroni & Cheese (frozen)"
$ws.Cells.Item(7, 22).Value = "'0.552"
$ws.Cells.Item(7, 22).Style = "Normal"
$ws.Cells.Item(7, 23).Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"
$ws.Cells.Item(7, 24).Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Cells.Item(7, 25).Value = "'0.664"
$ws.Cells.Item(7, 25).Style = "Normal"
$ws.Cells.Item(7, 26).Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"
$ws.Cells.Item(7, 27).Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Cells.Item(7, 28).Value = "'0.589"
$ws.Cells.Item(7, 28).Style = "Normal"
$ws.Cells.Item(7, 29).Value = "Portátil, saludable, fácil, buena textura, sabor suave"
$ws.Cells.Item(7, 30).Value = "Jack Link’s Beef Jerky Original"
$ws.Cells.Item(7, 31).Value = "'0.576"
$ws.Cells.Item(7, 31).Style = "Normal"
$ws.Cells.Item(7, 32).Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"

# --- Step 2: fix row 6 in place -- SmartScore columns become real numbers. ---
$ws.Cells.Item(6, 7).Value = 0.575
$ws.Cells.Item(6, 10).Value = 0.51
$ws.Cells.Item(6, 13).Value = 0.509
$ws.Cells.Item(6, 16).Value = 0.65
$ws.Cells.Item(6, 19).Value = 0.587
$ws.Cells.Item(6, 22).Value = 0.552
$ws.Cells.Item(6, 25).Value = 0.664
$ws.Cells.Item(6, 28).Value = 0.589
$ws.Cells.Item(6, 31).Value = 0.576

